{"js": "// Target edit: the last paragraph currently reads (split across two runs,\n// with a \"_GoBack\" bookmark sandwiched in between):\n//   run1: \"\u53bb\u6253\u9488\uff0c\u4e94\u74f6\uff01\u4e0d\u5982\u6740\u4e86\u6211\u3002\uff08\u65e0\u6240\u8c13\uff0c\u6211\u4f1a\u5b89\u8be6\u6253\u9488\"\n//   <bookmarkStart/><bookmarkEnd/>\n//   run2: \"\uff09\"\n//\n// After the edit it should read as a single run (bookmark untouched):\n//   run1: \"\u53bb\u6253\u9488\uff0c\u4e94\u74f6\uff01\u4e0d\u5982\u6740\u4e86\u6211\u3002\uff08\u65e0\u6240\u8c13\uff0c\u6211\u4f1a\u5b89\u8be6\u6253\u9488\uff09\uff0c\u65e0\u8bed\u3002\"\n//   <bookmarkStart/><bookmarkEnd/>\n// (run2 is gone \u2014 its \"\uff09\" text got folded into run1, plus the new \"\uff0c\u65e0\u8bed\u3002\")\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The edited paragraph is the last one in the document.\nconst target = paragraphs.items[paragraphs.items.length - 1];\n\n// Locate the trailing \"\uff09\" run's text (there is exactly one in this\n// paragraph) and delete it \u2014 this removes the whole second <w:r> and\n// leaves the bookmark sitting right after the first run.\nconst hits = target.search(\"\uff09\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[hits.items.length - 1].delete();\n  await context.sync();\n}\n\n// Append the new text at the end of the paragraph. Because the trailing\n// run was just removed, this lands right after run1 and merges into it\n// (inheriting run1's formatting), reproducing \"\uff09\uff0c\u65e0\u8bed\u3002\" being tacked\n// onto the original run's text.\ntarget.insertText(\"\uff09\uff0c\u65e0\u8bed\u3002\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Target edit: the last paragraph currently reads (split across two runs,\n# with a \"_GoBack\" bookmark sandwiched in between):\n#   run1: \"\u53bb\u6253\u9488\uff0c\u4e94\u74f6\uff01\u4e0d\u5982\u6740\u4e86\u6211\u3002\uff08\u65e0\u6240\u8c13\uff0c\u6211\u4f1a\u5b89\u8be6\u6253\u9488\"\n#   <bookmarkStart/><bookmarkEnd/>\n#   run2: \"\uff09\"\n#\n# After the edit it should read as a single run (bookmark untouched):\n#   run1: \"\u53bb\u6253\u9488\uff0c\u4e94\u74f6\uff01\u4e0d\u5982\u6740\u4e86\u6211\u3002\uff08\u65e0\u6240\u8c13\uff0c\u6211\u4f1a\u5b89\u8be6\u6253\u9488\uff09\uff0c\u65e0\u8bed\u3002\"\n#   <bookmarkStart/><bookmarkEnd/>\n# (run2 is gone -- its \"\uff09\" text got folded into run1, plus the new \"\uff0c\u65e0\u8bed\u3002\")\n\n$d = $word.ActiveDocument\n\n# Locate the edited paragraph by its distinctive content rather than assuming\n# a fixed index.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*\u53bb\u6253\u9488*\") {\n        $target = $p\n    }\n}\n\n$paraStart = $target.Range.Start\n$paraEnd = $target.Range.End\n\n# Find the LAST \"\uff09\" inside the paragraph (there happens to be just one, but\n# walk to the final hit to be safe against repeats).\n$lastFoundStart = -1\n$lastFoundEnd = -1\n$search = $d.Range($paraStart, $paraEnd)\nwhile ($search.Find.Execute(\"\uff09\")) {\n    if ($search.Start -lt $paraStart -or $search.End -gt $paraEnd) {\n        break\n    }\n    $lastFoundStart = $search.Start\n    $lastFoundEnd = $search.End\n    $search = $d.Range($search.End, $paraEnd)\n}\n\n# Delete that trailing \"\uff09\" run entirely -- this removes the whole second\n# <w:r> and leaves the bookmark sitting right after the first run.\nif ($lastFoundStart -ge 0) {\n    $toDelete = $d.Range($lastFoundStart, $lastFoundEnd)\n    $toDelete.Delete()\n}\n\n# Append the new text at the end of the paragraph. Because the trailing run\n# was just removed, this lands right after run1 and merges into it\n# (inheriting run1's formatting), reproducing \"\uff09\uff0c\u65e0\u8bed\u3002\" being tacked onto\n# the original run's text.\n$target.Range.InsertAfter(\"\uff09\uff0c\u65e0\u8bed\u3002\")\n"}
